$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 16: averages of the "biginteger" group (rows 6:8)
$ws.Range("B16").Formula = "=SUM(B6:B8)/3"
$ws.Range("C16").Formula = "=SUM(C6:C8)/3"
$ws.Range("D16:H16").Formula = "=SUM(D6:D8)/3"

# Row 17: averages of the "word" group (rows 9:11)
$ws.Range("C17").Formula = "=SUM(C9:C11)/3"
$ws.Range("D17:H17").Formula = "=SUM(D9:D11)/3"

# Row 18: averages of the next group (rows 12:14)
$ws.Range("C18").Formula = "=SUM(C12:C14)/3"
$ws.Range("D18:H18").Formula = "=SUM(D12:D14)/3"

# Update selection to match the author's final cursor position
$ws.Range("A20").Select()
